$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '29.237.61'
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').Value = '  +0.21%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '1.854.98'
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').Value = '  +0.14%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '0.7010'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +3.07%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '237.84'
$ws.Range('D6').Style = "Normal"
$ws.Range('E7').Value = '  -0.03%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.08057'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  +4.59%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.3022'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  -0.42%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '23.50'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  +1.77%  '
$ws.Range('E11').Value = '  +0.24%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '1.857.91'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  -0.29%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '5.200'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  +0.31%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '0.7059'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  -2.25%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '89.81'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  +0.66%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '29.282.36'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  +0.42%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '5.832'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  +2.02%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '0.000007878'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  +0.79%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '13.28'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  +1.22%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '236.97'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  +1.39%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '1.001'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  +0.01%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '2.110.99'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  +1.11%  '
$ws.Range('E23').Value = '  +0.01%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '7.468'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  +0.57%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '163.04'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  +0.74%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '8.886'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  -0.73%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '0.1408'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  -1.18%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '18.07'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  +0.21%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '1.917'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  -1.89%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '1.417'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  +1.57%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '1.473'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  -0.95%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '4.358'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  -3.50%  '
$ws.Range('E33').Value = '  +0.46%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.05190'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  +0.42%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '1.163'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  -1.31%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.7169'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  +2.22%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '1.000'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  -2.39%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '2.689'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  +0.77%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.01851'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  +0.44%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '2.726'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  +1.81%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.9317'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  +2.39%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '1.145.99'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  +4.44%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '6.020'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  +0.63%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.4253'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  -0.33%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '70.34'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  +0.39%  '
$ws.Range('E46').Value = '  -0.01%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '102.97'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  +0.61%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.5283'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  -3.61%  '
$ws.Range('B49').Value = 'RocketPoolETH'
$ws.Range('C49').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '2.008.52'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  +0.80%  '
$ws.Range('B50').Value = 'RenderToken'
$ws.Range('C50').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '1.744'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  -1.03%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '9.155'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  +0.33%  '
